# Add a new "Height" property row to the Property sheet for NPC.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

$newRow = 44

$ws.Cells.Item($newRow, 1).Value = "Height"
$ws.Cells.Item($newRow, 2).Value = "float"
$ws.Cells.Item($newRow, 3).Value = $false
$ws.Cells.Item($newRow, 4).Value = $false
$ws.Cells.Item($newRow, 5).Value = $false
$ws.Cells.Item($newRow, 6).Value = $true
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).NumberFormat = "@"
$ws.Cells.Item($newRow, 9).Value = "Friend"
$ws.Cells.Item($newRow, 10).Value = "模型高度"

# Match the author's final selection/viewport state (cosmetic, matches the
# commit's saved view: active cell on the newly-added row, scrolled down).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J" + $newRow).Select() | Out-Null
